$d = $word.ActiveDocument

# Locate the end of the last paragraph's text ("... tous mes coéquipiers."),
# right before the trailing _GoBack bookmark markers.
$r = $d.Content
$found = $r.Find.Execute("tous mes coéquipiers.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $splitPoint = $r.End

    # Split the paragraph at that point: this creates a new paragraph
    # (inheriting the same paragraph properties) that keeps the bookmark
    # markers, immediately after the existing "...coéquipiers." paragraph.
    $breakRange = $d.Range($splitPoint, $splitPoint)
    $breakRange.InsertBefore("`r")

    # The new (now last) paragraph starts right before the bookmark markers.
    $newPara = $d.Paragraphs.Last
    $insPos = $newPara.Range.Start

    $boldLabel = "Mardi 28 octobre 2014 : "
    $bodyText = "J’ai créé avec l’aide de Gabriel le diagramme de classe de nos contrôles visuels."
    $fullText = $boldLabel + $bodyText

    $textRange = $d.Range($insPos, $insPos)
    $textRange.InsertBefore($fullText)

    # Match the surrounding body text size (sz 24 half-points = 12pt).
    $wholeRange = $d.Range($insPos, $insPos + $fullText.Length)
    $wholeRange.Font.Size = 12

    # Bold only the date label.
    $boldRange = $d.Range($insPos, $insPos + $boldLabel.Length)
    $boldRange.Bold = 1
}
